# Update email / number values and the active selection on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix typo'd email address and corresponding number in row 2.
$ws.Range("F2").Value = "jonny.andrew@gmail.com"
$ws.Range("J2").Value = "788967494"

# Move the active selection from J2 to F4.
$ws.Range("F4").Select()
